$d = $word.ActiveDocument

$d.Content.Find.Execute("para la 2 convocatoria:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "para la 2ª y 3ª convocatoria:", 2)
